# Updated symbol list on Mon Jan 16 06:45:18 UTC 2023 with GitHub Actions
# Applies the refreshed Price / Volume(1h) figures (and, for the rows whose
# coin ranking shifted, the Coin name + Link) coming from the upstream feed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new text value. All of these columns are stored
# as plain text in the workbook (prices/percentages are text, not numbers),
# so we write them with a leading apostrophe to force text entry and then
# reset the style back to Normal so no stray number-format/quote-prefix
# style is left behind on the cell.
$updates = [ordered]@{
    'D2' = '304.04'
    'E2' = '2.09%'
    'D3' = '31.92'
    'E3' = '0.53%'
    'D4' = '5.223'
    'E4' = '2.19%'
    'D5' = '0.07846'
    'E5' = '4.29%'
    'D6' = '2.351'
    'E6' = '36.98%'
    'D7' = '7.999'
    'E7' = '3.14%'
    'B8' = 'MXToken'
    'C8' = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
    'D8' = '0.9125'
    'E8' = '-1.96%'
    'B9' = 'WazirX'
    'C9' = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
    'D9' = '0.1735'
    'E9' = '1.86%'
    'B10' = 'LiechtensteinCryptoassetsExchange'
    'C10' = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
    'D10' = '0.07387'
    'E10' = '-1.48%'
    'B11' = 'MandalaExchangeToken'
    'C11' = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
    'D11' = '0.08126'
    'E11' = '1.75%'
    'B12' = 'BitrueCoin'
    'C12' = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
    'D12' = '0.03049'
    'E12' = '0.43%'
    'B13' = 'BitMartToken'
    'C13' = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
    'D13' = '0.09946'
    'E13' = '0.50%'
    'B14' = 'BitForexToken'
    'C14' = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
    'D14' = '0.001515'
    'E14' = '1.73%'
    'B15' = 'TigerCash'
    'C15' = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
    'D15' = '0.006188'
    'E15' = '-3.87%'
    'B16' = 'LEO'
    'C16' = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
    'D16' = '3.500'
    'E16' = '1.07%'
    'B17' = 'GateToken'
    'C17' = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
    'D17' = '3.862'
    'E17' = '1.82%'
    'D18' = '2.239'
    'E18' = '0.83%'
    'D19' = '0.3298'
    'E19' = '0.54%'
    'D20' = '0.1311'
    'E20' = '-1.08%'
    'D21' = '4.660'
    'E21' = '2.07%'
    'D22' = '0.04636'
    'E22' = '-0.23%'
    'E23' = '0.41%'
    'E24' = '3.45%'
    'D25' = '0.004538'
    'E25' = '2.64%'
    'E26' = '3.74%'
    'D27' = '0.0002739'
    'E27' = '47.44%'
    'D39' = '0.01789'
    'E39' = '6.69%'
    'D40' = '0.04589'
    'E40' = '1.13%'
    'D41' = '0.007313'
    'E41' = '3.90%'
    'D42' = '0.1363'
    'E42' = '2.66%'
    'D43' = '0.002238'
    'E43' = '8.63%'
    'D44' = '0.01093'
    'E44' = '-7.33%'
    'D45' = '0.00006457'
    'E45' = '7.24%'
    'E46' = '-0.11%'
    'E47' = '15.31%'
    'D48' = '0.009890'
    'E48' = '-23.72%'
    'E49' = '-0.11%'
    'E50' = '-0.04%'
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    $cell.Value = "'" + $updates[$ref]
    $cell.Style = "Normal"
}

